# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed
# handback: status text, handback file/datetime columns, and the
# widened columns that the longer text now needs.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$mdFile  = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4c7c153358f6095858dd7d36bf6abcd5f8b4e2f9/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$zhXlf   = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.zh-cn.xlf"
$deXlf   = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.de-de.xlf"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: Status columns (E/F) for both rows now read the
# "handed back" message instead of "Ready for handoff".
# ---------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$ws1.Columns.Item(5).ColumnWidth = 29.1
$ws1.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------
# zh-cn sheet: Status (C) text changes via the shared string above;
# Latest Target File (I) / Latest Handback File (J) get filled in
# for both data rows, and the handback datetime (K) is refreshed.
# ---------------------------------------------------------------
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

$ws2.Range("I2").Value = $mdFile
$ws2.Range("J2").Value = $zhXlf
$ws2.Range("K2").Value = "2016-08-21 01:07:25"

$ws2.Range("I3").Value = $mdFile
$ws2.Range("J3").Value = $zhXlf
$ws2.Range("K3").Value = "2016-08-21 01:07:25"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("I3"), $mdUrl, "", "", $mdFile)

$ws2.Range("I2").Font.Underline = 2
$ws2.Range("I2").Font.Color = 15570276
$ws2.Range("I3").Font.Underline = 2
$ws2.Range("I3").Font.Color = 15570276

$ws2.Columns.Item(3).ColumnWidth = 29.1
$ws2.Columns.Item(9).ColumnWidth = 39.1667
$ws2.Columns.Item(10).ColumnWidth = 39.1667

# ---------------------------------------------------------------
# de-de sheet: same shape of edit, but the handback datetime is a
# brand-new timestamp (not reusing the zh-cn one).
# ---------------------------------------------------------------
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

$ws3.Range("I2").Value = $mdFile
$ws3.Range("J2").Value = $deXlf
$ws3.Range("K2").Value = "2016-08-21 01:07:31"

$ws3.Range("I3").Value = $mdFile
$ws3.Range("J3").Value = $deXlf
$ws3.Range("K3").Value = "2016-08-21 01:07:31"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("I3"), $mdUrl, "", "", $mdFile)

$ws3.Range("I2").Font.Underline = 2
$ws3.Range("I2").Font.Color = 15570276
$ws3.Range("I3").Font.Underline = 2
$ws3.Range("I3").Font.Color = 15570276

$ws3.Columns.Item(3).ColumnWidth = 29.1
$ws3.Columns.Item(9).ColumnWidth = 39.1667
$ws3.Columns.Item(10).ColumnWidth = 39.1667
